$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.609.50'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '1.666.84'
$ws.Range("E3").Value = '  -3.37%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.509'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.41%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.11'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.263'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = '1.904.23'
$ws.Range("E12").Value = '  -3.23%  '
$ws.Range("D13").Value = '1.677.50'
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("E14").Value = '  -3.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.560'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("D17").Value = '27.599.02'
$ws.Range("E17").Value = '  -1.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").Value = '0.0₃0731'
$ws.Range("E19").Value = '  -3.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.20%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("E24").Value = '  -3.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("E26").Value = '  -3.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("E29").Value = '  -2.17%  '
$ws.Range("E30").Value = '  +3.25%  '
$ws.Range("E31").Value = '  -1.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("D33").Value = '1.466.55'
$ws.Range("E33").Value = '  -2.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.49%  '
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.928'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.577'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.32%  '
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '69.66'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.03'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.82%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("E43").Value = '  -6.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.99%  '
$ws.Range("D45").Value = '1.810.50'
$ws.Range("E45").Value = '  -3.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.789'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("E49").Value = '  -5.27%  '
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.87%  '
